$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old data rows (2-10) - this also prunes now-unused shared strings
# (ECs, FAPs, MuSCs, Il12a, Il12rb2) so they can be re-added in the new order.
$ws.Range("A2:T10").Clear()

# Write the refreshed TPM data column-by-column (so the shared-string table
# is rebuilt in the same first-seen order as the target workbook: FAPs, MuSCs,
# Il12a, Il12rb2, ECs).

# Column A
$ws.Cells.Item(2,1).Value = "FAPs"
$ws.Cells.Item(3,1).Value = "FAPs"
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(5,1).Value = "MuSCs"
$ws.Cells.Item(6,1).Value = "MuSCs"
$ws.Cells.Item(7,1).Value = "MuSCs"

# Column B
$ws.Cells.Item(2,2).Value = "Il12a"
$ws.Cells.Item(3,2).Value = "Il12a"
$ws.Cells.Item(4,2).Value = "Il12a"
$ws.Cells.Item(5,2).Value = "Il12a"
$ws.Cells.Item(6,2).Value = "Il12a"
$ws.Cells.Item(7,2).Value = "Il12a"

# Column C
$ws.Cells.Item(2,3).Value = "Il12rb2"
$ws.Cells.Item(3,3).Value = "Il12rb2"
$ws.Cells.Item(4,3).Value = "Il12rb2"
$ws.Cells.Item(5,3).Value = "Il12rb2"
$ws.Cells.Item(6,3).Value = "Il12rb2"
$ws.Cells.Item(7,3).Value = "Il12rb2"

# Column D
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(7,4).Value = "MuSCs"

# Column E
$ws.Cells.Item(2,5).Value = 3.0
$ws.Cells.Item(3,5).Value = 3.0
$ws.Cells.Item(4,5).Value = 3.0
$ws.Cells.Item(5,5).Value = 2.0
$ws.Cells.Item(6,5).Value = 2.0
$ws.Cells.Item(7,5).Value = 2.0

# Column F
$ws.Cells.Item(2,6).Value = 1.0
$ws.Cells.Item(3,6).Value = 1.0
$ws.Cells.Item(4,6).Value = 1.0
$ws.Cells.Item(5,6).Value = 0.6666666666666666
$ws.Cells.Item(6,6).Value = 0.6666666666666666
$ws.Cells.Item(7,6).Value = 0.6666666666666666

# Column G
$ws.Cells.Item(2,7).Value = 2.338794666666666
$ws.Cells.Item(3,7).Value = 2.338794666666666
$ws.Cells.Item(4,7).Value = 2.338794666666666
$ws.Cells.Item(5,7).Value = 0.358666
$ws.Cells.Item(6,7).Value = 0.358666
$ws.Cells.Item(7,7).Value = 0.358666

# Column H
$ws.Cells.Item(2,8).Value = 7.016384
$ws.Cells.Item(3,8).Value = 7.016384
$ws.Cells.Item(4,8).Value = 7.016384
$ws.Cells.Item(5,8).Value = 1.075998
$ws.Cells.Item(6,8).Value = 1.075998
$ws.Cells.Item(7,8).Value = 1.075998

# Column I
$ws.Cells.Item(2,9).Value = 0.8670356886266615
$ws.Cells.Item(3,9).Value = 0.8670356886266615
$ws.Cells.Item(4,9).Value = 0.8670356886266615
$ws.Cells.Item(5,9).Value = 0.1329643113733386
$ws.Cells.Item(6,9).Value = 0.1329643113733386
$ws.Cells.Item(7,9).Value = 0.1329643113733386

# Column J
$ws.Cells.Item(2,10).Value = 0.8670356886266616
$ws.Cells.Item(3,10).Value = 0.8670356886266616
$ws.Cells.Item(4,10).Value = 0.8670356886266616
$ws.Cells.Item(5,10).Value = 0.1329643113733386
$ws.Cells.Item(6,10).Value = 0.1329643113733386
$ws.Cells.Item(7,10).Value = 0.1329643113733386

# Column K
$ws.Cells.Item(2,11).Value = 2.0
$ws.Cells.Item(3,11).Value = 3.0
$ws.Cells.Item(4,11).Value = 3.0
$ws.Cells.Item(5,11).Value = 2.0
$ws.Cells.Item(6,11).Value = 3.0
$ws.Cells.Item(7,11).Value = 3.0

# Column L
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(3,12).Value = 1.0
$ws.Cells.Item(4,12).Value = 1.0
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(6,12).Value = 1.0
$ws.Cells.Item(7,12).Value = 1.0

# Column M
$ws.Cells.Item(2,13).Value = 0.08759600000000001
$ws.Cells.Item(3,13).Value = 0.4379896666666667
$ws.Cells.Item(4,13).Value = 0.1255636666666667
$ws.Cells.Item(5,13).Value = 0.08759600000000001
$ws.Cells.Item(6,13).Value = 0.4379896666666667
$ws.Cells.Item(7,13).Value = 0.1255636666666667

# Column N
$ws.Cells.Item(2,14).Value = 0.262788
$ws.Cells.Item(3,14).Value = 1.313969
$ws.Cells.Item(4,14).Value = 0.376691
$ws.Cells.Item(5,14).Value = 0.262788
$ws.Cells.Item(6,14).Value = 1.313969
$ws.Cells.Item(7,14).Value = 0.376691

# Column O
$ws.Cells.Item(2,15).Value = 0.1345252087590763
$ws.Cells.Item(3,15).Value = 0.6726408893402845
$ws.Cells.Item(4,15).Value = 0.1928339019006393
$ws.Cells.Item(5,15).Value = 0.1345252087590763
$ws.Cells.Item(6,15).Value = 0.6726408893402845
$ws.Cells.Item(7,15).Value = 0.1928339019006393

# Column P
$ws.Cells.Item(2,16).Value = 0.1345252087590763
$ws.Cells.Item(3,16).Value = 0.6726408893402844
$ws.Cells.Item(4,16).Value = 0.1928339019006393
$ws.Cells.Item(5,16).Value = 0.1345252087590763
$ws.Cells.Item(6,16).Value = 0.6726408893402844
$ws.Cells.Item(7,16).Value = 0.1928339019006393

# Column Q
$ws.Cells.Item(2,17).Value = 0.2048690576213333
$ws.Cells.Item(3,17).Value = 1.024367896455111
$ws.Cells.Item(4,17).Value = 0.293667633927111
$ws.Cells.Item(5,17).Value = 0.031417706936
$ws.Cells.Item(6,17).Value = 0.1570920017846667
$ws.Cells.Item(7,17).Value = 0.04503541806866666

# Column R
$ws.Cells.Item(2,18).Value = 1.843821518592
$ws.Cells.Item(3,18).Value = 9.219311068096001
$ws.Cells.Item(4,18).Value = 2.643008705344
$ws.Cells.Item(5,18).Value = 0.282759362424
$ws.Cells.Item(6,18).Value = 1.413828016062
$ws.Cells.Item(7,18).Value = 0.405318762618

# Column S
$ws.Cells.Item(2,19).Value = 0.1166381570140711
$ws.Cells.Item(3,19).Value = 0.5832036566876037
$ws.Cells.Item(4,19).Value = 0.1671938749249869
$ws.Cells.Item(5,19).Value = 0.01788705174500519
$ws.Cells.Item(6,19).Value = 0.08943723265268096
$ws.Cells.Item(7,19).Value = 0.02564002697565242

# Column T
$ws.Cells.Item(2,20).Value = 0.1166381570140711
$ws.Cells.Item(3,20).Value = 0.5832036566876035
$ws.Cells.Item(4,20).Value = 0.1671938749249869
$ws.Cells.Item(5,20).Value = 0.01788705174500519
$ws.Cells.Item(6,20).Value = 0.08943723265268094
$ws.Cells.Item(7,20).Value = 0.02564002697565242
